# Update the "PACS Submit status" cell (Q2) on Sheet1 from "Successfully "
# to "Successfully Verifired".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("Q2").Value = "Successfully Verifired"
